# "Added sra ids for under_ice_rerun"
#
# BF3:BF25  -> Bioproject accession PRJNAXXXXXX is replaced, in place, by the
#              real accession PRJNA417044 (shared by every row).
# BG3:BG25  -> Biosample placeholder is replaced by the real accessions,
#              one per physical sample (rows split into three groups).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("micans_v6_ex1")

# Update the Bioproject accession everywhere it appears (rewrites the shared
# string in place instead of minting a brand-new one).
$ws.Cells.Replace("PRJNAXXXXXX", "PRJNA417044")

# Normalize style (the refreshed cells in the source workbook lost their
# explicit style index and fall back to the default cell style).
$ws.Range("BF3:BF25").Style = "Normal"

# Biosample accessions, grouped by physical sample.
$ws.Range("BG3:BG10").Value = "SAMN07975454"
$ws.Range("BG11:BG18").Value = "SAMN07975455"
$ws.Range("BG19:BG25").Value = "SAMN07975456"
$ws.Range("BG3:BG25").Style = "Normal"

# Restore the selection left behind on the sheet (bottom pane scrolled to the
# last edited cell).
$ws.Range("BG29").Select()
